$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add four new timesheet entries in rows 29-32 (two shifts on 11/2/2024,
#    two shifts on 11/4/2024), matching the formatting of similarly-styled
#    existing rows so the cell styles line up with the rest of the sheet.
# ---------------------------------------------------------------------------

# Row 29 -> style donor row 19 (date / h:mm / h:mm AM/PM)
$ws.Range("A19:G19").Copy() | Out-Null
$ws.Range("A29:G29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 30 -> style donor row 17 (date / h:mm AM/PM / h:mm)
$ws.Range("A17:G17").Copy() | Out-Null
$ws.Range("A30:G30").PasteSpecial(-4122) | Out-Null

# Row 31 -> style donor row 28 (date / h:mm AM/PM / h:mm AM/PM)
$ws.Range("A28:G28").Copy() | Out-Null
$ws.Range("A31:G31").PasteSpecial(-4122) | Out-Null

# Row 32 -> style donor row 28 (date / h:mm AM/PM / h:mm AM/PM)
$ws.Range("A28:G28").Copy() | Out-Null
$ws.Range("A32:G32").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Date + time values for the new entries
$ws.Range("A29").Value = 45598
$ws.Range("B29").Value = 0.46875
$ws.Range("C29").Value = 0.54166666666666663

$ws.Range("A30").Value = 45598
$ws.Range("B30").Value = 0.58750000000000002
$ws.Range("C30").Value = 0.61805555555555558

$ws.Range("A31").Value = 45600
$ws.Range("B31").Value = 0.41736111111111113
$ws.Range("C31").Value = 0.59375

$ws.Range("A32").Value = 45600
$ws.Range("B32").Value = 0.81597222222222221
$ws.Range("C32").Value = 0.87291666666666667

# Fill the D:G derived-value formulas down across the new rows as one shared
# formula block (same pattern already used for D5:D28 / E5:E28 / F5:F28 /
# G3:G28).
$ws.Range("D29:D32").Formula = "=C29-B29"
$ws.Range("E29:E32").Formula = "=D29*1440"
$ws.Range("F29:F32").Formula = "=E29/60"
$ws.Range("G29:G32").Formula = "=F29*22.5"

# ---------------------------------------------------------------------------
# 2. New weekly-total cells for the just-added rows.
# ---------------------------------------------------------------------------
$ws.Range("M6:N6").Copy() | Out-Null
$ws.Range("M7:N7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("M7").Formula = "=SUM(D29:D32)"
$ws.Range("N7").Formula = "=SUM(G29:G32)"

# ---------------------------------------------------------------------------
# 3. Clean up the old leftover stub row (D31) and move the
#    "DON'T FORGET TO SAVE (local)" reminder from row 33 down to row 36 now
#    that real data occupies rows 29-32.
# ---------------------------------------------------------------------------
$ws.Range("B33").Delete(-4162) | Out-Null   # xlShiftUp

$ws.Range("B36").Value = "DON'T FORGET TO SAVE (local)"
$ws.Range("B36").Font.Bold = $true

# Match the author's final cursor position.
$ws.Range("L19").Select() | Out-Null

$wb.Save()
